$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values are plain text such as "1.000" or "29.366.70"
# rather than numbers, so a text number format is applied before the
# assignment to stop Excel from re-interpreting them as numeric/date
# values; the cell style is then restored to the workbook default so no
# extra formatting is introduced.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Refresh price (D) and 1h volume change (E) for each listed coin ---
Set-TextValue $ws.Cells.Item(2, 4) "29.366.70"
$ws.Cells.Item(2, 5).Value = "  -0.14%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.840.72"
$ws.Cells.Item(3, 5).Value = "  -0.41%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9988"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
Set-TextValue $ws.Cells.Item(5, 4) "239.37"
$ws.Cells.Item(5, 5).Value = "  -0.41%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.6264"
$ws.Cells.Item(6, 5).Value = "  +0.05%  "
Set-TextValue $ws.Cells.Item(7, 4) "1.000"
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.07438"
$ws.Cells.Item(8, 5).Value = "  -0.79%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.2899"
$ws.Cells.Item(9, 5).Value = "  -0.14%  "
Set-TextValue $ws.Cells.Item(10, 4) "24.81"
$ws.Cells.Item(10, 5).Value = "  +1.57%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07718"
$ws.Cells.Item(11, 5).Value = "  -0.28%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.835.89"
$ws.Cells.Item(12, 5).Value = "  -0.68%  "
Set-TextValue $ws.Cells.Item(13, 4) "4.966"
$ws.Cells.Item(13, 5).Value = "  -0.73%  "
Set-TextValue $ws.Cells.Item(14, 4) "0.6762"
$ws.Cells.Item(14, 5).Value = "  -0.67%  "
Set-TextValue $ws.Cells.Item(15, 4) "0.00001020"
$ws.Cells.Item(15, 5).Value = "  -2.43%  "
Set-TextValue $ws.Cells.Item(16, 4) "81.74"
$ws.Cells.Item(16, 5).Value = "  -0.62%  "
Set-TextValue $ws.Cells.Item(17, 4) "6.230"
$ws.Cells.Item(17, 5).Value = "  +0.94%  "
Set-TextValue $ws.Cells.Item(18, 4) "29.340.18"
$ws.Cells.Item(18, 5).Value = "  -0.33%  "
Set-TextValue $ws.Cells.Item(19, 4) "231.87"
$ws.Cells.Item(19, 5).Value = "  +0.88%  "
Set-TextValue $ws.Cells.Item(20, 4) "12.30"
$ws.Cells.Item(20, 5).Value = "  -0.39%  "
Set-TextValue $ws.Cells.Item(21, 4) "1.000"
$ws.Cells.Item(21, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(22, 4) "7.360"
$ws.Cells.Item(22, 5).Value = "  -1.55%  "
Set-TextValue $ws.Cells.Item(23, 4) "1.001"
$ws.Cells.Item(23, 5).Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(24, 4) "158.38"
$ws.Cells.Item(24, 5).Value = "  -0.47%  "
Set-TextValue $ws.Cells.Item(25, 4) "8.481"
$ws.Cells.Item(25, 5).Value = "  +0.80%  "
Set-TextValue $ws.Cells.Item(26, 4) "0.1349"
$ws.Cells.Item(26, 5).Value = "  -1.76%  "
Set-TextValue $ws.Cells.Item(27, 4) "17.34"
$ws.Cells.Item(27, 5).Value = "  -1.20%  "
Set-TextValue $ws.Cells.Item(28, 4) "0.07224"
$ws.Cells.Item(28, 5).Value = "  +13.00%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.459"
$ws.Cells.Item(29, 5).Value = "  +4.26%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.481"
$ws.Cells.Item(30, 5).Value = "  +0.25%  "
Set-TextValue $ws.Cells.Item(31, 4) "4.045"
$ws.Cells.Item(31, 5).Value = "  -1.20%  "
Set-TextValue $ws.Cells.Item(32, 4) "4.043"
$ws.Cells.Item(32, 5).Value = "  -1.25%  "
Set-TextValue $ws.Cells.Item(33, 4) "1.818"
$ws.Cells.Item(33, 5).Value = "  -0.59%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.140"
$ws.Cells.Item(34, 5).Value = "  -0.18%  "
Set-TextValue $ws.Cells.Item(35, 4) "0.6977"
$ws.Cells.Item(35, 5).Value = "  -0.13%  "
Set-TextValue $ws.Cells.Item(36, 4) "2.573"
$ws.Cells.Item(36, 5).Value = "  -0.18%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.01842"
$ws.Cells.Item(37, 5).Value = "  +0.68%  "
Set-TextValue $ws.Cells.Item(38, 4) "6.941"
$ws.Cells.Item(38, 5).Value = "  +4.71%  "
Set-TextValue $ws.Cells.Item(39, 4) "2.815"
$ws.Cells.Item(39, 5).Value = "  -0.35%  "
Set-TextValue $ws.Cells.Item(40, 4) "1.234.46"
$ws.Cells.Item(40, 5).Value = "  -2.43%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.9408"
$ws.Cells.Item(41, 5).Value = "  +3.37%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.9999"
$ws.Cells.Item(42, 5).Value = "  +0.08%  "
Set-TextValue $ws.Cells.Item(43, 4) "1.988.93"
$ws.Cells.Item(43, 5).Value = "  -1.01%  "
Set-TextValue $ws.Cells.Item(44, 4) "100.67"
$ws.Cells.Item(44, 5).Value = "  -0.80%  "
Set-TextValue $ws.Cells.Item(45, 4) "65.67"
$ws.Cells.Item(45, 5).Value = "  -0.96%  "
Set-TextValue $ws.Cells.Item(48, 4) "6.953"
$ws.Cells.Item(48, 5).Value = "  -1.74%  "
Set-TextValue $ws.Cells.Item(49, 4) "8.902"
$ws.Cells.Item(49, 5).Value = "  -1.56%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.1138"
$ws.Cells.Item(50, 5).Value = "  -2.96%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.3903"
$ws.Cells.Item(51, 5).Value = "  -1.32%  "

# --- RenderToken and BabyDogeCoin swapped ranking positions (rows 46/47) ---
$ws.Cells.Item(46, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Cells.Item(46, 4) "0.00000000119"
$ws.Cells.Item(46, 5).Value = "  +1.94%  "

$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Cells.Item(47, 4) "1.718"
$ws.Cells.Item(47, 5).Value = "  -1.57%  "
